$d = $word.ActiveDocument

$find = "constel" + [char]0xB7 + "lació, "
$replace = " "

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replace, 2)
